$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells to match new layout.
# F1: JEFE INMEDIATO -> CEDULA DEL JEFE INMEDIATO
$ws.Range("F1").Value = "CÉDULA DEL JEFE INMEDIATO"
# I1: TURNO -> HORA DE ENTRADA L-V
$ws.Range("I1").Value = "HORA DE ENTRADA L-V"
# J1: ESTADO -> HORA DE SALIDA L-V
$ws.Range("J1").Value = "HORA DE SALIDA L-V"
# New columns K, L (write L before K so shared-string order matches source)
$ws.Range("L1").Value = "HORA DE SALIDA SÁBADOS"
$ws.Range("K1").Value = "HORA DE ENTRADA SÁBADOS"

# Copy the header style from an existing header cell onto the two new cells
# so they reuse the same cell-style index instead of generating new ones.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Column widths for the new "hora" columns (I and J hold the longer labels).
$ws.Columns.Item(9).ColumnWidth = 20.43
$ws.Columns.Item(10).ColumnWidth = 18.6

# Header row is taller to accommodate wrapped, longer labels.
$ws.Rows.Item(1).RowHeight = 45

# Move the active selection to the new last header cell.
$ws.Range("L1").Select() | Out-Null
